$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "MENDEZ"
$ws.Range("C6").Value = "ESQUITE"
$ws.Range("E6").Value = "ANDY"
$ws.Range("G6").Value = "ELIAS"
$ws.Range("I6").Value = "/201762626"

$ws.Range("A12").Value = "2009-03-31"
$ws.Range("F12").Value = "8A"

$ws.Range("A14").Value = ""
$ws.Range("F14").Value = ""

$ws.Range("A18").Value = ""
$ws.Range("F18").Value = "HEIDY MENDEZ"

$ws.Range("A24").Value = "24/10/2017"
$ws.Range("C24").Value = "14:48:24"
